# Rename the `name` attribute (wp:docPr/@name and pic:cNvPr/@name) that Word
# stores for three inline pictures: the Pearson logo in both footers and the
# BTEC logo in the "first page" header. Word's object model has no writable
# InlineShape.Name (only Shape.Name for floating shapes), so the rename is
# done by round-tripping just that picture run's WordOpenXML through
# Range.InsertXML, which is the supported way to edit markup that isn't
# reachable through a scalar OM property.

function Rename-InlineShapeMarkupName {
    param(
        $Shape,
        [string]$OldName,
        [string]$NewName
    )

    $rng = $Shape.Range
    $xml = $rng.WordOpenXML

    # Isolate just the <w:r>...</w:r> run that owns the <w:drawing>, so the
    # InsertXML round-trip only has to re-mint that one run instead of the
    # whole header/footer story.
    $drawIdx = $xml.IndexOf("<w:drawing>")
    $runStart = $xml.LastIndexOf("<w:r ", $drawIdx)
    $runEndIdx = $xml.IndexOf("</w:r>", $drawIdx) + "</w:r>".Length
    $runXml = $xml.Substring($runStart, $runEndIdx - $runStart)

    $newRunXml = $runXml.Replace('name="' + $OldName + '"', 'name="' + $NewName + '"')

    $bodyStart = $xml.IndexOf("<w:body>") + "<w:body>".Length
    $bodyEnd = $xml.IndexOf("</w:body>")

    $prefix = $xml.Substring(0, $bodyStart)
    $suffix = $xml.Substring($bodyEnd)

    $newXml = $prefix + $newRunXml + $suffix
    $rng.InsertXML($newXml)
}

$d = $word.ActiveDocument
$sec = $d.Sections.First

# Footer "default" (footer2.xml): Pearson logo docPr id="2" -> image2.png to image1.png
$footerDefault = $sec.Footers.Item(1)
Rename-InlineShapeMarkupName -Shape $footerDefault.Range.InlineShapes.Item(1) -OldName "image2.png" -NewName "image1.png"

# Footer "first page" (footer1.xml): Pearson logo docPr id="3" -> image2.png to image1.png
$footerFirst = $sec.Footers.Item(2)
Rename-InlineShapeMarkupName -Shape $footerFirst.Range.InlineShapes.Item(1) -OldName "image2.png" -NewName "image1.png"

# Header "first page" (header1.xml): BTEC logo docPr id="1" -> image1.jpg to image2.jpg
$headerFirst = $sec.Headers.Item(2)
Rename-InlineShapeMarkupName -Shape $headerFirst.Range.InlineShapes.Item(1) -OldName "image1.jpg" -NewName "image2.jpg"
